$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("A5").Value = 'L474792'
$ws.Range("C5").Value = 'NICK SHELL SERVICE'
$ws.Range("E5").Value = 1840
$ws.Range("F5").Value = 'T'
$ws.Range("H5").Value = 45241.0421827199
$ws.Range("I5").ClearContents()
$ws.Range("J5").Value = '10/15/23 13:31'
$ws.Range("K5").Value = '10/15/23 13:31'
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = '$1,840 as of 10/15/2023 11:31:29 AM'
$ws.Range("N5").Value = 1880
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 0
$ws.Range("Q5").Value = 0

# Row 6
$ws.Range("A6").Value = 'L647934'
$ws.Range("C6").Value = 'SB #6'
$ws.Range("E6").Value = 1940
$ws.Range("F6").Value = 'T'
$ws.Range("H6").ClearContents()
$ws.Range("I6").Value = 'ATM Inactive greater than 2000 minutes'
$ws.Range("J6").Value = '04/06/23 22:10'
$ws.Range("K6").Value = '04/06/23 22:05'
$ws.Range("L6").Value = 20
$ws.Range("M6").Value = '$1,940 as of 4/6/2023 8:05:45 PM'
$ws.Range("N6").Value = 1960
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = 0
$ws.Range("Q6").Value = 0

# Row 7
$ws.Range("A7").Value = 'LK644532'
$ws.Range("C7").Value = 'SCL ENTERPRISES LAUNDRY'
$ws.Range("E7").Value = 2320
$ws.Range("F7").Value = 'T'
$ws.Range("H7").Value = 45283.0421827199
$ws.Range("I7").ClearContents()
$ws.Range("J7").Value = '10/16/23 15:08'
$ws.Range("K7").Value = '10/16/23 15:08'
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = '$2,340 as of 10/15/2023 8:14:28 PM'
$ws.Range("N7").Value = 2340
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = 0
$ws.Range("Q7").Value = 0

# Row 8
$ws.Range("A8").Value = 'L678988'
$ws.Range("C8").Value = 'PAYELESS MARKET'
$ws.Range("E8").Value = 2400
$ws.Range("F8").Value = 'T'
$ws.Range("H8").ClearContents()
$ws.Range("I8").Value = 'ATM Inactive greater than 2000 minutes'
$ws.Range("J8").Value = '07/20/23 20:09'
$ws.Range("K8").Value = '07/20/23 20:09'
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = '$2,400 as of 7/20/2023 6:09:40 PM'
$ws.Range("N8").Value = 2500
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = 0
$ws.Range("Q8").Value = 0

# Row 9
$ws.Range("A9").Value = 'L688966'
$ws.Range("C9").Value = 'S B WESTERN 108TH MARKET'
$ws.Range("E9").Value = 2580
$ws.Range("F9").Value = 'T'
$ws.Range("H9").Value = 45224.0421827199
$ws.Range("I9").ClearContents()
$ws.Range("J9").Value = '10/16/23 18:43'
$ws.Range("K9").Value = '10/16/23 18:43'
$ws.Range("L9").Value = 80
$ws.Range("M9").Value = '$2,600 as of 10/16/2023 10:06:02 AM'
$ws.Range("N9").Value = 2600
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = 0
$ws.Range("Q9").Value = 0

# Row 10
$ws.Range("A10").Value = 'LK561655'
$ws.Range("C10").Value = 'CRENSHAW CRAVOR #2'
$ws.Range("E10").Value = 2780
$ws.Range("F10").Value = 'T'
$ws.Range("H10").ClearContents()
$ws.Range("I10").Value = 'ATM Inactive greater than 48 minutes'
$ws.Range("J10").Value = '01/23/20 08:24'
$ws.Range("K10").Value = '01/23/20 08:24'
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = '$2,780 as of 1/23/2020 6:24:32 AM'
$ws.Range("N10").Value = 2800
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = 0
$ws.Range("Q10").Value = 0

# Row 11
$ws.Range("A11").Value = 'L474761'
$ws.Range("C11").Value = 'BABS MARKET'
$ws.Range("E11").Value = 3560
$ws.Range("F11").Value = 'T'
$ws.Range("H11").Value = 45279.0421827199
$ws.Range("I11").ClearContents()
$ws.Range("J11").Value = '10/16/23 19:07'
$ws.Range("K11").Value = '10/16/23 19:07'
$ws.Range("L11").Value = 100
$ws.Range("M11").Value = '$3,660 as of 10/14/2023 1:54:54 PM'
$ws.Range("N11").Value = 3660
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = 0
$ws.Range("Q11").Value = 0

# Row 12
$ws.Range("A12").Value = 'L475182'
$ws.Range("C12").Value = 'LA ESQUINA DE ORO'
$ws.Range("E12").Value = 3800
$ws.Range("F12").Value = 'T'
$ws.Range("H12").ClearContents()
$ws.Range("I12").Value = 'ATM Inactive greater than 48 minutes'
$ws.Range("J12").Value = '09/16/20 16:57'
$ws.Range("K12").Value = '09/15/20 23:38'
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = '$3,800 as of 9/16/2020 1:28:00 PM'
$ws.Range("N12").Value = 3800
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = 0
$ws.Range("Q12").Value = 0

# Row 13
$ws.Range("A13").Value = 'L488595'
$ws.Range("C13").Value = 'N S MART'
$ws.Range("E13").Value = 4060
$ws.Range("F13").Value = 'T'
$ws.Range("H13").Value = 45272.0421827199
$ws.Range("I13").ClearContents()
$ws.Range("J13").Value = '10/15/23 22:41'
$ws.Range("K13").Value = '10/15/23 15:07'
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = '$4,060 as of 10/15/2023 8:41:09 PM'
$ws.Range("N13").Value = 4060
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = 0
$ws.Range("Q13").Value = 0

# Row 14
$ws.Range("A14").Value = 'L662336'
$ws.Range("C14").Value = 'SB#4 MONA MARKET'
$ws.Range("E14").Value = 4580
$ws.Range("F14").Value = 'T'
$ws.Range("H14").Value = 45232.0421827199
$ws.Range("I14").ClearContents()
$ws.Range("J14").Value = '10/15/23 14:42'
$ws.Range("K14").Value = '10/15/23 14:42'
$ws.Range("L14").Value = 120
$ws.Range("M14").Value = '$4,580 as of 10/15/2023 12:42:04 PM'
$ws.Range("N14").Value = 4580
$ws.Range("O14").Value = 0
$ws.Range("P14").Value = 0
$ws.Range("Q14").Value = 0

# Row 15
$ws.Range("A15").Value = 'L697590'
$ws.Range("C15").Value = 'S B MARKET ST'
$ws.Range("E15").Value = 5400
$ws.Range("F15").Value = 'T'
$ws.Range("H15").Value = 45274.0421827199
$ws.Range("I15").ClearContents()
$ws.Range("J15").Value = '10/16/23 15:58'
$ws.Range("K15").Value = '10/16/23 15:58'
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = '$5,600 as of 10/15/2023 7:53:24 PM'
$ws.Range("N15").Value = 5600
$ws.Range("O15").Value = 0
$ws.Range("P15").Value = 0
$ws.Range("Q15").Value = 0

# Row 16
$ws.Range("A16").Value = 'LK864765'
$ws.Range("C16").Value = 'SKY LIQUOR'
$ws.Range("E16").Value = 5560
$ws.Range("F16").Value = 'T'
$ws.Range("H16").Value = 45230.0421827199
$ws.Range("I16").ClearContents()
$ws.Range("J16").Value = '10/16/23 18:47'
$ws.Range("K16").Value = '10/16/23 15:18'
$ws.Range("L16").Value = 60
$ws.Range("M16").Value = '$5,720 as of 10/16/2023 1:38:43 AM'
$ws.Range("N16").Value = 5620
$ws.Range("O16").Value = 0
$ws.Range("P16").Value = 0
$ws.Range("Q16").Value = 0

# Row 17
$ws.Range("A17").Value = 'LK236828'
$ws.Range("C17").Value = 'WORLDWIDE AUTOMOTIVE'
$ws.Range("E17").Value = 5760
$ws.Range("F17").Value = 'T'
$ws.Range("H17").Value = 45241.0421827199
$ws.Range("I17").ClearContents()
$ws.Range("J17").Value = '10/16/23 17:49'
$ws.Range("K17").Value = '10/16/23 17:49'
$ws.Range("L17").Value = 20
$ws.Range("M17").Value = '$6,100 as of 10/14/2023 5:10:56 PM'
$ws.Range("N17").Value = 5860
$ws.Range("O17").Value = 0
$ws.Range("P17").Value = 0
$ws.Range("Q17").Value = 0

# Row 18
$ws.Range("A18").Value = 'L476340'
$ws.Range("C18").Value = 'DONUT & SANDWICH'
$ws.Range("E18").Value = 5800
$ws.Range("F18").Value = 'T'
$ws.Range("H18").Value = 45237.0421827199
$ws.Range("I18").ClearContents()
$ws.Range("J18").Value = '10/16/23 15:35'
$ws.Range("K18").Value = '10/16/23 11:14'
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = '$5,800 as of 10/16/2023 9:14:56 AM'
$ws.Range("N18").Value = 5800
$ws.Range("O18").Value = 0
$ws.Range("P18").Value = 0
$ws.Range("Q18").Value = 0

# Row 19
$ws.Range("A19").Value = 'L704741'
$ws.Range("C19").Value = 'W ADAMS COIN LAUNDRY'
$ws.Range("E19").Value = 6140
$ws.Range("F19").Value = 'T'
$ws.Range("H19").Value = 45222.0421827199
$ws.Range("I19").ClearContents()
$ws.Range("J19").Value = '10/16/23 17:46'
$ws.Range("K19").Value = '10/16/23 17:46'
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = '$6,300 as of 10/16/2023 11:08:11 AM'
$ws.Range("N19").Value = 6160
$ws.Range("O19").Value = 0
$ws.Range("P19").Value = 0
$ws.Range("Q19").Value = 0

# Row 20
$ws.Range("A20").Value = 'L474817'
$ws.Range("C20").Value = 'SAFETY MARKET'
$ws.Range("E20").Value = 6660
$ws.Range("F20").Value = 'T'
$ws.Range("H20").Value = 45229.0421827199
$ws.Range("I20").ClearContents()
$ws.Range("J20").Value = '10/16/23 18:28'
$ws.Range("K20").Value = '10/16/23 18:28'
$ws.Range("L20").Value = 100
$ws.Range("M20").Value = '$6,720 as of 10/16/2023 10:06:34 AM'
$ws.Range("N20").Value = 6660
$ws.Range("O20").Value = 0
$ws.Range("P20").Value = 0
$ws.Range("Q20").Value = 0

# Row 21
$ws.Range("A21").Value = 'L682801'
$ws.Range("C21").Value = 'SB#5'
$ws.Range("E21").Value = 7840
$ws.Range("F21").Value = 'T'
$ws.Range("H21").ClearContents()
$ws.Range("I21").Value = 'ATM Inactive greater than 2000 minutes'
$ws.Range("J21").Value = '09/28/23 15:22'
$ws.Range("K21").Value = '09/28/23 12:14'
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = '$7,840 as of 9/28/2023 12:31:50 PM'
$ws.Range("N21").Value = 7840
$ws.Range("O21").Value = 0
$ws.Range("P21").Value = 0
$ws.Range("Q21").Value = 0

# Row 22
$ws.Range("A22").Value = 'L474746'
$ws.Range("C22").Value = 'ZACATES MARKET'
$ws.Range("E22").Value = 8020
$ws.Range("F22").Value = 'T'
$ws.Range("H22").Value = 45262.0421827199
$ws.Range("I22").ClearContents()
$ws.Range("J22").Value = '10/16/23 14:19'
$ws.Range("K22").Value = '10/16/23 14:19'
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = '$8,080 as of 10/16/2023 11:09:43 AM'
$ws.Range("N22").Value = 8080
$ws.Range("O22").Value = 0
$ws.Range("P22").Value = 0
$ws.Range("Q22").Value = 0

# Row 23
$ws.Range("A23").Value = 'L475090'
$ws.Range("C23").Value = 'S.B. 2'
$ws.Range("E23").Value = 9320
$ws.Range("F23").Value = 'T'
$ws.Range("H23").Value = 45239.0421827199
$ws.Range("I23").ClearContents()
$ws.Range("J23").Value = '10/16/23 13:40'
$ws.Range("K23").Value = '10/16/23 13:40'
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = '$9,320 as of 10/16/2023 11:40:18 AM'
$ws.Range("N23").Value = 9400
$ws.Range("O23").Value = 0
$ws.Range("P23").Value = 0
$ws.Range("Q23").Value = 0

# Row 24
$ws.Range("A24").Value = 'L688961'
$ws.Range("C24").Value = 'MONA MART'
$ws.Range("E24").Value = 9480
$ws.Range("F24").Value = 'T'
$ws.Range("H24").Value = 46605.0421827199
$ws.Range("I24").ClearContents()
$ws.Range("J24").Value = '10/16/23 18:29'
$ws.Range("K24").Value = '10/16/23 16:18'
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = '$9,500 as of 10/16/2023 10:30:50 AM'
$ws.Range("N24").Value = 9480
$ws.Range("O24").Value = 0
$ws.Range("P24").Value = 0
$ws.Range("Q24").Value = 0

# Row 25
$ws.Range("A25").Value = 'L697589'
$ws.Range("C25").Value = 'S B DISCOUNT MART'
$ws.Range("E25").Value = 12300
$ws.Range("F25").Value = 'T'
$ws.Range("H25").Value = 45228.0421827199
$ws.Range("I25").ClearContents()
$ws.Range("J25").Value = '10/16/23 19:00'
$ws.Range("K25").Value = '10/16/23 19:00'
$ws.Range("L25").Value = 60
$ws.Range("M25").Value = '$12,460 as of 10/16/2023 11:54:58 AM'
$ws.Range("N25").Value = 12460
$ws.Range("O25").Value = 0
$ws.Range("P25").Value = 0
$ws.Range("Q25").Value = 0

# Row 26
$ws.Range("A26").Value = 'LK923383'
$ws.Range("C26").Value = 'SAMYS PHONE CARDS'
$ws.Range("E26").Value = 12320
$ws.Range("F26").Value = 'T'
$ws.Range("H26").Value = 45237.0421827199
$ws.Range("I26").ClearContents()
$ws.Range("J26").Value = '10/16/23 17:53'
$ws.Range("K26").Value = '10/16/23 17:53'
$ws.Range("L26").Value = 80
$ws.Range("M26").Value = '$12,640 as of 10/16/2023 10:57:51 AM'
$ws.Range("N26").Value = 12520
$ws.Range("O26").Value = 0
$ws.Range("P26").Value = 0
$ws.Range("Q26").Value = 0

# Row 27
$ws.Range("A27").Value = 'LK891176'
$ws.Range("C27").Value = '98 DISCOUNT STORE'
$ws.Range("E27").Value = 21040
$ws.Range("F27").Value = 'T'
$ws.Range("H27").Value = 45232.0421827199
$ws.Range("I27").ClearContents()
$ws.Range("J27").Value = '10/16/23 19:04'
$ws.Range("K27").Value = '10/16/23 18:59'
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = '$21,080 as of 10/16/2023 9:46:55 AM'
$ws.Range("N27").Value = 21040
$ws.Range("O27").Value = 0
$ws.Range("P27").Value = 0
$ws.Range("Q27").Value = 0

# Row 28 - Total
$ws.Range("E28").Value = 145500